# chore: update Sheets via scheduled runner
# Refreshes cached Universalis market-price snapshots (currentAveragePrice /
# currentAveragePriceNQ/HQ) and the derived Leve profit columns for the rows
# whose item prices moved since the last run.
$wb = $excel.ActiveWorkbook


# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 12
$ws.Range("H12").Value = 1095.25
$ws.Range("I12").Value = 861.4
$ws.Range("K12").Value = 861.4
$ws.Range("M12").Value = -691.4
# row 33
$ws.Range("H33").Value = 186.46666
$ws.Range("I33").Value = 55.88889
$ws.Range("K33").Value = 55.88889
$ws.Range("M33").Value = 173.11111
# row 74
$ws.Range("H74").Value = 5489.9
$ws.Range("I74").Value = 5552.6665
$ws.Range("J74").Value = 5395.75
$ws.Range("K74").Value = 5552.6665
$ws.Range("L74").Value = 5395.75
$ws.Range("M74").Value = -4616.6665
$ws.Range("N74").Value = -7267.75
# row 77
$ws.Range("H77").Value = 5489.9
$ws.Range("I77").Value = 5552.6665
$ws.Range("J77").Value = 5395.75
$ws.Range("K77").Value = 27763.3325
$ws.Range("L77").Value = 26978.75
$ws.Range("M77").Value = -23083.3325
$ws.Range("N77").Value = -36338.75
# row 92
$ws.Range("H92").Value = 901.6429000000001
$ws.Range("I92").Value = 718.6667
$ws.Range("J92").Value = 1999.5
$ws.Range("K92").Value = 718.6667
$ws.Range("L92").Value = 1999.5
$ws.Range("M92").Value = 529.3333
$ws.Range("N92").Value = -4495.5
# row 96
$ws.Range("H96").Value = 2969.875
$ws.Range("I96").Value = 1770.3334
$ws.Range("K96").Value = 5311.0002
$ws.Range("M96").Value = -3938.0002

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 7268.256
$ws.Range("I32").Value = 5449.615
$ws.Range("K32").Value = 5449.615
$ws.Range("M32").Value = -5162.615
# row 97
$ws.Range("H97").Value = 1083.1666
$ws.Range("J97").Value = 966.6667
$ws.Range("L97").Value = 966.6667
$ws.Range("N97").Value = -1958.6667
# row 102
$ws.Range("H102").Value = 1902.5
$ws.Range("I102").Value = 1055
$ws.Range("K102").Value = 1055
$ws.Range("M102").Value = 567
# row 110
$ws.Range("H110").Value = 2766.9167
$ws.Range("I110").Value = 3603.4443
$ws.Range("K110").Value = 3603.4443
$ws.Range("M110").Value = -1558.4443
# row 132
$ws.Range("H132").Value = 1576.0952
$ws.Range("I132").Value = 1481.9744
$ws.Range("K132").Value = 4445.9232
$ws.Range("M132").Value = -1915.9232

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 99
$ws.Range("H99").Value = 4560.3335
$ws.Range("I99").Value = 4453.9
$ws.Range("J99").Value = 4773.2
$ws.Range("K99").Value = 4453.9
$ws.Range("L99").Value = 4773.2
$ws.Range("M99").Value = -2955.9
$ws.Range("N99").Value = -7769.2

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 55
$ws.Range("H55").Value = 5000
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 5000
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 5000
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -5630
# row 68
$ws.Range("H68").Value = 32499.5
$ws.Range("J68").Value = 32499.5
$ws.Range("L68").Value = 32499.5
$ws.Range("N68").Value = -33997.5
# row 71
$ws.Range("H71").Value = 32499.5
$ws.Range("J71").Value = 32499.5
$ws.Range("L71").Value = 97498.5
$ws.Range("N71").Value = -104986.5
# row 74
$ws.Range("H74").Value = 47782
$ws.Range("J74").Value = 45314
$ws.Range("L74").Value = 45314
$ws.Range("N74").Value = -47062
# row 77
$ws.Range("H77").Value = 47782
$ws.Range("J77").Value = 45314
$ws.Range("L77").Value = 135942
$ws.Range("N77").Value = -144678
# row 132
$ws.Range("H132").Value = 2425.0476
$ws.Range("I132").Value = 1596.9231
$ws.Range("K132").Value = 4790.7693
$ws.Range("M132").Value = -2260.7693

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 124
$ws.Range("H124").Value = 1049.5
$ws.Range("I124").Value = 1049.5
$ws.Range("K124").Value = 3148.5
$ws.Range("M124").Value = 1761.5
# row 125
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()
# row 131
$ws.Range("H131").Value = 1867.25
$ws.Range("J131").Value = 1891.3334
$ws.Range("L131").Value = 5674.0002
$ws.Range("N131").Value = -15754.0002

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 132
$ws.Range("H132").Value = 2032.0952
$ws.Range("I132").Value = 1432.5555
$ws.Range("J132").Value = 5629.3335
$ws.Range("K132").Value = 4297.666499999999
$ws.Range("L132").Value = 16888.0005
$ws.Range("M132").Value = -1767.666499999999
$ws.Range("N132").Value = -21948.0005

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Range("H22").Value = 9266.4
$ws.Range("I22").Value = 4533.8
$ws.Range("J22").Value = 13999
$ws.Range("K22").Value = 4533.8
$ws.Range("L22").Value = 13999
$ws.Range("M22").Value = -4238.8
$ws.Range("N22").Value = -14589
# row 27
$ws.Range("H27").Value = 9266.4
$ws.Range("I27").Value = 4533.8
$ws.Range("J27").Value = 13999
$ws.Range("K27").Value = 4533.8
$ws.Range("L27").Value = 13999
$ws.Range("M27").Value = -4426.8
$ws.Range("N27").Value = -14213
# row 40
$ws.Range("H40").Value = 1973.25
$ws.Range("I40").Value = 1973.25
$ws.Range("K40").Value = 1973.25
$ws.Range("M40").Value = -1837.25
# row 61
$ws.Range("H61").Value = 4710.8335
$ws.Range("I61").Value = 5413.4443
$ws.Range("J61").Value = 4008.2222
$ws.Range("K61").Value = 5413.4443
$ws.Range("L61").Value = 4008.2222
$ws.Range("M61").Value = -5211.4443
$ws.Range("N61").Value = -4412.2222
# row 93
$ws.Range("H93").Value = 300
$ws.Range("I93").Value = 305
$ws.Range("K93").Value = 305
$ws.Range("M93").Value = 943
# row 100
$ws.Range("H100").Value = 3680.3635
$ws.Range("J100").Value = 4655
$ws.Range("L100").Value = 4655
$ws.Range("N100").Value = -5737
# row 113
$ws.Range("H113").Value = 4710.8335
$ws.Range("I113").Value = 5413.4443
$ws.Range("J113").Value = 4008.2222
$ws.Range("K113").Value = 5413.4443
$ws.Range("L113").Value = 4008.2222
$ws.Range("M113").Value = -3243.4443
$ws.Range("N113").Value = -8348.2222

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 69
$ws.Range("H69").Value = 20000
$ws.Range("I69").Value = 20000
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 20000
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -19251
$ws.Range("N69").ClearContents()
# row 72
$ws.Range("H72").Value = 20000
$ws.Range("I72").Value = 20000
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 60000
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -56256
$ws.Range("N72").ClearContents()
